$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.285.16'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.39%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.902.09'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.36%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '326.42'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.40%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.25%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4652'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3927'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07890'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -0.93%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9893'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '21.99'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.73%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.901.63'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -3.92%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.081'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.05%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.751'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -0.91%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.06977'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +0.22%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '88.39'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.47%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.002'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.12%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000009988'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '17.08'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -0.76%  '
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.23%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '29.283.55'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +0.23%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.317'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.91%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.10'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.096'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.61%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '156.15'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '19.41'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.61%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '5.987'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +2.01%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '118.63'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.85%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.907'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -5.06%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09356'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -0.68%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.9074'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -1.76%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.285'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -1.52%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.328'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.49%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.215'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -1.51%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.187'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.77%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.05784'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02092'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.66%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '7.741'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -3.65%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5714'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -0.87%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1785'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.61%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '9.764'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -2.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.91'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.62%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5352'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.47%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.191'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.90%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.07040'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.80%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -1.70%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.569'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.33%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '113.60'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +0.98%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.059'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.49%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '71.22'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.82%  '
